$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 34, column A with the refined timestamp value
$ws.Cells.Item(34, 1).Value = 44347.8294425625

# Append new row 35 with the newly retrieved data
$ws.Cells.Item(35, 1).Value = 44348.86069004962
$ws.Cells.Item(35, 2).Value = 74631
$ws.Cells.Item(35, 3).Value = 62928
$ws.Cells.Item(35, 4).Value = 3172
$ws.Cells.Item(35, 5).Value = 2072
$ws.Cells.Item(35, 6).Value = 1472
$ws.Cells.Item(35, 7).Value = 19523
$ws.Cells.Item(35, 8).Value = 1339
$ws.Cells.Item(35, 9).Value = 880
$ws.Cells.Item(35, 10).Value = 203
